$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sales rows to append (producto, precio, unidades, metodo_pago, fecha, total)
$rows = @(
    @("Samsung Galaxy S21",     15000, 1, "efectivo", "2025-06-16", 15000),
    @("Xiaomi Redmi Note 12",    7000, 1, "efectivo", "2025-06-16",  7000),
    @("Motorola G73",            6500, 3, "tarjeta",  "2025-06-16", 19500),
    @("Samsung Galaxy S21",     15000, 1, "tarjeta",  "2025-06-16", 15000),
    @("Samsung Galaxy S21",     15000, 2, "tarjeta",  "2025-06-16", 30000),
    @("Motorola G73",            6500, 1, "tarjeta",  "2025-06-16",  6500),
    @("Xiaomi Redmi Note 12",    7000, 2, "tarjeta",  "2025-06-16", 14000),
    @("Samsung Galaxy S21",     15000, 1, "efectivo", "2025-06-16", 15000),
    @("Motorola G73",            6500, 2, "efectivo", "2025-06-16", 13000),
    @("Xiaomi Redmi Note 12",    7000, 1, "tarjeta",  "2025-06-16",  7000),
    @("Samsung Galaxy S21",     15000, 3, "efectivo", "2025-06-16", 45000),
    @("iPhone 13",              18000, 3, "tarjeta",  "2025-06-16", 54000)
)

$r = 5
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Force the date column to remain plain text (matches source data which
    # stores "fecha" as literal text, not a date serial), then clear the
    # temporary text format so no stray style sticks to the cell.
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).ClearFormats()

    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
